# Update the cover-page version/date string in preparation for the
# 11.10.01 release:
#   "Version 11.08.01, 2016-02-15"  ->  "Version 11.10.01, 2016-05-02"
#
# Only four digits actually change:
#   ".08."  -> ".10."   (the "0" becomes "1", the "8" becomes "0")
#   "-02-"  -> "-05-"   (the "2" becomes "5")
#   "-15"   -> "-02"    (the day "15" becomes "02")

$d = $word.ActiveDocument

# Locate the exact string robustly (rather than hard-coding absolute
# character offsets) so the edit still lands correctly even if content
# earlier in the document shifts.
$anchor = $d.Content.Duplicate
$anchor.Find.ClearFormatting()
$anchor.Find.Execute("Version 11.08.01, 2016-02-15", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$base = $anchor.Start

# Apply the four single-spot edits from the right (end of the string)
# back to the left. Editing right-to-left keeps every not-yet-visited
# offset stable (no downstream shifting to account for) and keeps the
# untouched text before each edit point undisturbed.

# "15" -> "02"  (day)
$d.Range($base + 26, $base + 28).Text = "02"

# "2" -> "5"  (month, inside "-02-")
$d.Range($base + 24, $base + 25).Text = "5"

# "8" -> "0"  (inside ".08.")
$d.Range($base + 12, $base + 13).Text = "0"

# "0" -> "1"  (inside ".08.")
$d.Range($base + 11, $base + 12).Text = "1"
